# Scheduled-runner market data refresh: update currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ columns (H-N) for a set of leve rows across
# each crafting-job sheet with freshly pulled market values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3862.4546
$ws.Range("I40").Value = 5416.6665
$ws.Range("J40").Value = 1997.4
$ws.Range("K40").Value = 5416.6665
$ws.Range("L40").Value = 1997.4
$ws.Range("M40").Value = -5241.6665
$ws.Range("N40").Value = -2347.4

$ws.Range("H92").Value = 3199.4
$ws.Range("I92").Value = 3499.7222
$ws.Range("K92").Value = 3499.7222
$ws.Range("M92").Value = -2251.7222

$ws.Range("H101").Value = 1033.4445
$ws.Range("I101").Value = 685
$ws.Range("J101").Value = 2775.6667
$ws.Range("K101").Value = 2055
$ws.Range("L101").Value = 8327.000100000001
$ws.Range("M101").Value = -433
$ws.Range("N101").Value = -11571.0001

$ws.Range("H112").Value = 2523.5667
$ws.Range("J112").Value = 2888.28
$ws.Range("L112").Value = 8664.84
$ws.Range("N112").Value = -10880.84

$ws.Range("H125").Value = 3897.5715
$ws.Range("I125").Value = 4043.5386
$ws.Range("K125").Value = 36391.8474
$ws.Range("M125").Value = -33931.8474

$ws.Range("H137").Value = 9330.352000000001
$ws.Range("I137").Value = 1581
$ws.Range("J137").Value = 27735.062
$ws.Range("K137").Value = 4743
$ws.Range("L137").Value = 83205.186
$ws.Range("M137").Value = -2193
$ws.Range("N137").Value = -88305.186

$ws.Range("H138").Value = 3496.305
$ws.Range("J138").Value = 3307.3804
$ws.Range("L138").Value = 9922.1412
$ws.Range("N138").Value = -20202.1412

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 12484.576
$ws.Range("I61").Value = 2916.0476
$ws.Range("K61").Value = 2916.0476
$ws.Range("M61").Value = -2704.0476

$ws.Range("H74").Value = 16372.655
$ws.Range("I74").Value = 2539.8235
$ws.Range("J74").Value = 35969.168
$ws.Range("K74").Value = 2539.8235
$ws.Range("L74").Value = 35969.168
$ws.Range("M74").Value = -1665.8235
$ws.Range("N74").Value = -37717.168

$ws.Range("H77").Value = 16372.655
$ws.Range("I77").Value = 2539.8235
$ws.Range("J77").Value = 35969.168
$ws.Range("K77").Value = 12699.1175
$ws.Range("L77").Value = 179845.84
$ws.Range("M77").Value = -8331.1175
$ws.Range("N77").Value = -188581.84

$ws.Range("H122").Value = 1638475.1
$ws.Range("I122").Value = 3679065.5
$ws.Range("J122").Value = 6002.7
$ws.Range("K122").Value = 11037196.5
$ws.Range("L122").Value = 18008.1
$ws.Range("M122").Value = -11034746.5
$ws.Range("N122").Value = -22908.1

$ws.Range("H132").Value = 3586219
$ws.Range("I132").Value = 5609.6924
$ws.Range("J132").Value = 6689413.5
$ws.Range("K132").Value = 16829.0772
$ws.Range("L132").Value = 20068240.5
$ws.Range("M132").Value = -14299.0772
$ws.Range("N132").Value = -20073300.5

$ws.Range("H136").Value = 12484.576
$ws.Range("I136").Value = 2916.0476
$ws.Range("K136").Value = 8748.1428
$ws.Range("M136").Value = -6198.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 31811.637
$ws.Range("I96").Value = 16321.333
$ws.Range("K96").Value = 16321.333
$ws.Range("M96").Value = -13575.333

$ws.Range("H107").Value = 4945
$ws.Range("I107").Value = 6049.375
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 6049.375
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -4129.375
$ws.Range("N107").Value = -5840

$ws.Range("H134").Value = 15446.479
$ws.Range("I134").Value = 8064
$ws.Range("J134").Value = 29288.625
$ws.Range("K134").Value = 24192
$ws.Range("L134").Value = 87865.875
$ws.Range("M134").Value = -21657
$ws.Range("N134").Value = -92935.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 70939.05499999999
$ws.Range("I31").Value = 118716.445
$ws.Range("J31").Value = 20351.234
$ws.Range("K31").Value = 118716.445
$ws.Range("L31").Value = 20351.234
$ws.Range("M31").Value = -118421.445
$ws.Range("N31").Value = -20941.234

$ws.Range("H34").Value = 70939.05499999999
$ws.Range("I34").Value = 118716.445
$ws.Range("J34").Value = 20351.234
$ws.Range("K34").Value = 118716.445
$ws.Range("L34").Value = 20351.234
$ws.Range("M34").Value = -118514.445
$ws.Range("N34").Value = -20755.234

$ws.Range("H55").Value = 33374.75
$ws.Range("J55").Value = 33374.75
$ws.Range("L55").Value = 33374.75
$ws.Range("N55").Value = -34004.75

$ws.Range("H58").Value = 20063.84
$ws.Range("I58").Value = 7730.8184
$ws.Range("K58").Value = 7730.8184
$ws.Range("M58").Value = -7527.8184

$ws.Range("H99").Value = 4730731
$ws.Range("I99").Value = 4454400
$ws.Range("J99").Value = 5007062.5
$ws.Range("K99").Value = 4454400
$ws.Range("L99").Value = 5007062.5
$ws.Range("M99").Value = -4452902
$ws.Range("N99").Value = -5010058.5

$ws.Range("H107").Value = 1155.8108
$ws.Range("I107").Value = 886.62964
$ws.Range("K107").Value = 886.62964
$ws.Range("M107").Value = 1033.37036

$ws.Range("H126").Value = 4730731
$ws.Range("I126").Value = 4454400
$ws.Range("J126").Value = 5007062.5
$ws.Range("K126").Value = 13363200
$ws.Range("L126").Value = 15021187.5
$ws.Range("M126").Value = -13360730
$ws.Range("N126").Value = -15026127.5

$ws.Range("H132").Value = 5268.484
$ws.Range("I132").Value = 2255.2173
$ws.Range("K132").Value = 6765.651899999999
$ws.Range("M132").Value = -4235.651899999999

$ws.Range("H136").Value = 20063.84
$ws.Range("I136").Value = 7730.8184
$ws.Range("K136").Value = 23192.4552
$ws.Range("M136").Value = -20642.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 2342.8096
$ws.Range("I97").Value = 278.57144
$ws.Range("J97").Value = 6471.2856
$ws.Range("K97").Value = 835.71432
$ws.Range("L97").Value = 19413.8568
$ws.Range("M97").Value = -339.71432
$ws.Range("N97").Value = -20405.8568

$ws.Range("H129").Value = 5054890.5
$ws.Range("J129").Value = 6997850.5
$ws.Range("L129").Value = 20993551.5
$ws.Range("N129").Value = -21003551.5

$ws.Range("H131").Value = 1480.78
$ws.Range("J131").Value = 1480.78
$ws.Range("L131").Value = 4442.34
$ws.Range("N131").Value = -14522.34

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 15148.6
$ws.Range("J92").Value = 16054
$ws.Range("L92").Value = 16054
$ws.Range("N92").Value = -19798

$ws.Range("H107").Value = 1083.2222
$ws.Range("J107").Value = 1199.8572
$ws.Range("L107").Value = 1199.8572
$ws.Range("N107").Value = -5039.8572

$ws.Range("H122").Value = 1998319
$ws.Range("I122").Value = 2264161.5
$ws.Range("J122").Value = 4499.5
$ws.Range("K122").Value = 6792484.5
$ws.Range("L122").Value = 13498.5
$ws.Range("M122").Value = -6790034.5
$ws.Range("N122").Value = -18398.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1839831.4
$ws.Range("I7").Value = 2984476
$ws.Range("J7").Value = 8400
$ws.Range("K7").Value = 2984476
$ws.Range("L7").Value = 8400
$ws.Range("M7").Value = -2984364
$ws.Range("N7").Value = -8624

$ws.Range("H46").Value = 2767.0908
$ws.Range("I46").Value = 1370.5555
$ws.Range("K46").Value = 1370.5555
$ws.Range("M46").Value = -1182.5555

$ws.Range("H76").Value = 23014.5
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 23014.5
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 23014.5
$ws.Range("N76").Value = -23690.5
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 23014.5
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 23014.5
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 23014.5
$ws.Range("N79").Value = -25354.5
$ws.Range("M79").ClearContents()

$ws.Range("H126").Value = 1839831.4
$ws.Range("I126").Value = 2984476
$ws.Range("J126").Value = 8400
$ws.Range("K126").Value = 8953428
$ws.Range("L126").Value = 25200
$ws.Range("M126").Value = -8950958
$ws.Range("N126").Value = -30140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H126").Value = 6342.231
$ws.Range("I126").Value = 6342.231
$ws.Range("K126").Value = 19026.693
$ws.Range("M126").Value = -16556.693

$ws.Range("H128").Value = 85000
$ws.Range("J128").Value = 85000
$ws.Range("L128").Value = 85000
$ws.Range("N128").Value = -94960

$ws.Range("H129").Value = 22916.666
$ws.Range("J129").Value = 22916.666
$ws.Range("L129").Value = 22916.666
$ws.Range("N129").Value = -32916.666

$ws.Range("H132").Value = 4848.102
$ws.Range("I132").Value = 2253.2856
$ws.Range("K132").Value = 6759.8568
$ws.Range("M132").Value = -4229.8568
